$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.175.17"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.87%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.102.26"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +5.24%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.20"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.11"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +7.47%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.098.47"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +5.27%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.527"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.49%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.70"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.07%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +3.29%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.482"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +5.41%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000251"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +3.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.01"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +8.98%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.41%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.620.83"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.188.77"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.52%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.24"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +4.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.106.31"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +5.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.26"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +17.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "471.88"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +5.71%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.718"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +6.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.55"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +5.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.00"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.61%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.37"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +9.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.92"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +6.97%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.27"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +3.79%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.12"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +3.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.43"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +4.66%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +4.81%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0000102"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +5.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "28.46"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +4.38%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +5.65%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +4.48%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.94"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +4.06%  "
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "Arweave"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "47.35"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +10.88%  "
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.11"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +6.85%  "
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.47"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.87%  "
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "TheGraph"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.319"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +7.16%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +4.63%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.93"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +5.95%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.74"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +3.65%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "397.01"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.40%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0365"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +3.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.772.11"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.61%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "135.46"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +3.66%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.82"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +7.58%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.25"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +5.17%  "
